# Auto-generated Excel COM-interop script to apply the Behemoth_Profits data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 569.8333
$ws.Range("I12").Value = 429.75
$ws.Range("K12").Value = 429.75
$ws.Range("M12").Value = -259.75
$ws.Range("H96").Value = 2768.5715
$ws.Range("I96").Value = 1783.2222
$ws.Range("J96").Value = 4542.2
$ws.Range("K96").Value = 5349.6666
$ws.Range("L96").Value = 13626.6
$ws.Range("M96").Value = -3976.6666
$ws.Range("N96").Value = -16372.6
$ws.Range("H99").Value = 552.4
$ws.Range("I99").Value = 572.3333
$ws.Range("J99").Value = 522.5
$ws.Range("K99").Value = 1716.9999
$ws.Range("L99").Value = 1567.5
$ws.Range("M99").Value = -218.9999
$ws.Range("N99").Value = -4563.5
$ws.Range("H106").Value = 9367.267
$ws.Range("I106").Value = 1612.0
$ws.Range("K106").Value = 1612.0
$ws.Range("M106").Value = -981.0
$ws.Range("H138").Value = 3025.725
$ws.Range("I138").Value = 2255.2856
$ws.Range("J138").Value = 3099.6028
$ws.Range("K138").Value = 6765.8568
$ws.Range("L138").Value = 9298.8084
$ws.Range("M138").Value = -1625.8568
$ws.Range("N138").Value = -19578.8084

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1821.3334
$ws.Range("J21").Value = 2357.0
$ws.Range("L21").Value = 2357.0
$ws.Range("N21").Value = -3105.0
$ws.Range("H32").Value = 7359552.5
$ws.Range("I32").Value = 10002327.0
$ws.Range("J32").Value = 18513.777
$ws.Range("K32").Value = 10002327.0
$ws.Range("L32").Value = 18513.777
$ws.Range("M32").Value = -10002040.0
$ws.Range("N32").Value = -19087.777
$ws.Range("H46").Value = 3327.8333
$ws.Range("J46").Value = 3164.4
$ws.Range("L46").Value = 3164.4
$ws.Range("N46").Value = -3802.4
$ws.Range("H76").Value = 49999.332
$ws.Range("J76").Value = 49999.332
$ws.Range("L76").Value = 49999.332
$ws.Range("N76").Value = -50675.332
$ws.Range("H79").Value = 49999.332
$ws.Range("J79").Value = 49999.332
$ws.Range("L79").Value = 49999.332
$ws.Range("N79").Value = -52339.332
$ws.Range("H102").Value = 15783.6
$ws.Range("I102").Value = 15783.6
$ws.Range("J102").Value = 0.0
$ws.Range("K102").Value = 15783.6
$ws.Range("L102").Value = 0.0
$ws.Range("M102").Value = -14161.6
$ws.Range("N102").Value = $null
$ws.Range("H110").Value = 1623.125
$ws.Range("I110").Value = 1495.75
$ws.Range("J110").Value = 1750.5
$ws.Range("K110").Value = 1495.75
$ws.Range("L110").Value = 1750.5
$ws.Range("M110").Value = 549.25
$ws.Range("N110").Value = -5840.5
$ws.Range("H123").Value = 85000.0
$ws.Range("J123").Value = 85000.0
$ws.Range("L123").Value = 85000.0
$ws.Range("N123").Value = -94800.0
$ws.Range("H131").Value = 73329.0
$ws.Range("J131").Value = 73329.0
$ws.Range("L131").Value = 73329.0
$ws.Range("N131").Value = -83409.0
$ws.Range("H134").Value = 103678.94
$ws.Range("J134").Value = 103678.94
$ws.Range("L134").Value = 103678.94
$ws.Range("N134").Value = -113818.94

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8228.368
$ws.Range("I99").Value = 12846.333
$ws.Range("J99").Value = 4072.2
$ws.Range("K99").Value = 12846.333
$ws.Range("L99").Value = 4072.2
$ws.Range("M99").Value = -11348.333
$ws.Range("N99").Value = -7068.2
$ws.Range("H105").Value = 2269.8948
$ws.Range("I105").Value = 1535.7333
$ws.Range("J105").Value = 2748.6956
$ws.Range("K105").Value = 1535.7333
$ws.Range("L105").Value = 2748.6956
$ws.Range("M105").Value = 211.2666999999999
$ws.Range("N105").Value = -6242.6956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1119.5
$ws.Range("I7").Value = 248.25
$ws.Range("J7").Value = 2862.0
$ws.Range("K7").Value = 248.25
$ws.Range("L7").Value = 2862.0
$ws.Range("M7").Value = -135.25
$ws.Range("N7").Value = -3088.0
$ws.Range("H69").Value = 52921.0
$ws.Range("I69").Value = 24000.0
$ws.Range("J69").Value = 62561.332
$ws.Range("K69").Value = 24000.0
$ws.Range("L69").Value = 62561.332
$ws.Range("M69").Value = -23251.0
$ws.Range("N69").Value = -64059.332
$ws.Range("H72").Value = 52921.0
$ws.Range("I72").Value = 24000.0
$ws.Range("J72").Value = 62561.332
$ws.Range("K72").Value = 72000.0
$ws.Range("L72").Value = 187683.996
$ws.Range("M72").Value = -68256.0
$ws.Range("N72").Value = -195171.996
$ws.Range("H103").Value = 34445.145
$ws.Range("I103").Value = 14999.75
$ws.Range("J103").Value = 60372.332
$ws.Range("K103").Value = 14999.75
$ws.Range("L103").Value = 60372.332
$ws.Range("M103").Value = -13827.75
$ws.Range("N103").Value = -62716.332
$ws.Range("H108").Value = 81473.336
$ws.Range("J108").Value = 87310.0
$ws.Range("L108").Value = 87310.0
$ws.Range("N108").Value = -94990.0
$ws.Range("H129").Value = 65666.336
$ws.Range("J129").Value = 65666.336
$ws.Range("L129").Value = 65666.336
$ws.Range("N129").Value = -75666.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15975104.0
$ws.Range("I4").Value = 11255509.0
$ws.Range("K4").Value = 33766527.0
$ws.Range("M4").Value = -33766415.0
$ws.Range("H37").Value = 76052.555
$ws.Range("J37").Value = 76052.555
$ws.Range("L37").Value = 228157.665
$ws.Range("N37").Value = -228381.665
$ws.Range("H58").Value = 1761.1538
$ws.Range("J58").Value = 2785.7144
$ws.Range("L58").Value = 8357.143199999999
$ws.Range("N58").Value = -8613.143199999999
$ws.Range("H107").Value = 1678.0
$ws.Range("J107").Value = 1678.0
$ws.Range("L107").Value = 5034.0
$ws.Range("N107").Value = -8874.0
$ws.Range("H113").Value = 1461.0625
$ws.Range("I113").Value = 766.3333
$ws.Range("J113").Value = 1877.9
$ws.Range("K113").Value = 2298.9999
$ws.Range("L113").Value = 5633.700000000001
$ws.Range("M113").Value = -128.9998999999998
$ws.Range("N113").Value = -9973.7
$ws.Range("H117").Value = 2045.0
$ws.Range("J117").Value = 2045.0
$ws.Range("L117").Value = 6135.0
$ws.Range("N117").Value = -13019.0
$ws.Range("H121").Value = 816.9048
$ws.Range("J121").Value = 1195.5714
$ws.Range("L121").Value = 3586.7142
$ws.Range("N121").Value = -6206.7142
$ws.Range("H131").Value = 9384.947
$ws.Range("J131").Value = 10826.546
$ws.Range("L131").Value = 32479.638
$ws.Range("N131").Value = -42559.638
$ws.Range("H134").Value = 3946.1904

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 322.0
$ws.Range("I2").Value = 272.8889
$ws.Range("J2").Value = 385.14285
$ws.Range("K2").Value = 272.8889
$ws.Range("L2").Value = 385.14285
$ws.Range("M2").Value = -159.8889
$ws.Range("N2").Value = -611.14285
$ws.Range("H39").Value = 0.0
$ws.Range("J39").Value = 0.0
$ws.Range("L39").Value = 0.0
$ws.Range("N39").Value = $null
$ws.Range("H102").Value = 2352.9355
$ws.Range("J102").Value = 3543.1428
$ws.Range("L102").Value = 3543.1428
$ws.Range("N102").Value = -6787.1428
$ws.Range("H132").Value = 30305938.0
$ws.Range("I132").Value = 34485756.0
$ws.Range("J132").Value = 2250.0
$ws.Range("K132").Value = 103457268.0
$ws.Range("L132").Value = 6750.0
$ws.Range("M132").Value = -103454738.0
$ws.Range("N132").Value = -11810.0
$ws.Range("H136").Value = 6398.45
$ws.Range("J136").Value = 6398.45
$ws.Range("L136").Value = 19195.35
$ws.Range("N136").Value = -24295.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3616.3333
$ws.Range("I46").Value = 2939.4
$ws.Range("J46").Value = 4462.5
$ws.Range("K46").Value = 2939.4
$ws.Range("L46").Value = 4462.5
$ws.Range("M46").Value = -2751.4
$ws.Range("N46").Value = -4838.5
$ws.Range("H82").Value = 3359.8
$ws.Range("I82").Value = 0.0
$ws.Range("K82").Value = 0.0
$ws.Range("M82").Value = $null
$ws.Range("H85").Value = 3359.8
$ws.Range("I85").Value = 0.0
$ws.Range("K85").Value = 0.0
$ws.Range("M85").Value = $null
$ws.Range("H100").Value = 3334.1667
$ws.Range("I100").Value = 1002.5
$ws.Range("K100").Value = 1002.5
$ws.Range("M100").Value = -461.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7236.56
$ws.Range("I81").Value = 3705.8
$ws.Range("K81").Value = 7411.6
$ws.Range("M81").Value = -6350.6
$ws.Range("H84").Value = 7236.56
$ws.Range("I84").Value = 3705.8
$ws.Range("K84").Value = 37058.0
$ws.Range("M84").Value = -31754.0
$ws.Range("H96").Value = 2527.6365
$ws.Range("I96").Value = 2163.3333
$ws.Range("J96").Value = 2964.8
$ws.Range("K96").Value = 2163.3333
$ws.Range("L96").Value = 2964.8
$ws.Range("M96").Value = -790.3332999999998
$ws.Range("N96").Value = -5710.8
$ws.Range("H102").Value = 93730.0
$ws.Range("J102").Value = 93730.0
$ws.Range("L102").Value = 93730.0
$ws.Range("N102").Value = -100220.0
$ws.Range("H122").Value = 2604.625
$ws.Range("I122").Value = 2646.258
$ws.Range("J122").Value = 1314.0
$ws.Range("K122").Value = 7938.773999999999
$ws.Range("L122").Value = 3942.0
$ws.Range("M122").Value = -5488.773999999999
$ws.Range("N122").Value = -8842.0
$ws.Range("H132").Value = 254052.97
$ws.Range("I132").Value = 3594.6
$ws.Range("K132").Value = 10783.8
$ws.Range("M132").Value = -8253.8
$ws.Range("H136").Value = 3215.1282
$ws.Range("I136").Value = 2738.9429
$ws.Range("K136").Value = 8216.8287
$ws.Range("M136").Value = -5666.8287
